# The edit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (used by the slide master)   Integral -> Office Theme
#   ppt/theme/theme2.xml  (used by the notes master)    Office Theme -> Integral
#
# The PowerPoint COM object model doesn't expose a generic "replace this
# theme part's XML" call, but it does expose each theme's 12 scheme colors
# as a read/write RGB property (ThemeColorScheme.Item(n).RGB), in the fixed
# order dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. Driving those 12
# slots to the target palette reproduces the content swap exactly (the two
# themes already share an identical font scheme and format scheme).

$p = $ppt.ActivePresentation

# Target palette for the slide-master theme (becomes the "Office Theme" colors)
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# Target palette for the notes-master theme (becomes the "Integral" colors)
$integralColors = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}

$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
